$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point value in A12 (recalculated timestamp)
$ws.Range("A12").Value = 45876.41687603009

# Add new row 13 with the latest sensor reading
$ws.Range("A13").Value = 45876.4585312047
$ws.Range("B13").Value = 2025
$ws.Range("C13").Value = 28
$ws.Range("D13").Value = 18.56
$ws.Range("E13").Value = 81.63
$ws.Range("F13").Value = 580.08
$ws.Range("G13").Value = 11.91
$ws.Range("H13").Value = "ESE"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "11:00:17"

# Apply the same date/time number format as A12 to A13
$ws.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
